$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2504.3643659222066
$ws.Range("B1").Value = 1661.9251245798198
$ws.Range("C1").Value = 1646.5851732613664
$ws.Range("A2").Value = 2270.1765981343892
$ws.Range("B2").Value = 1579.8198104143721
$ws.Range("C2").Value = 1489.133339545851
$ws.Range("A3").Value = 2562.8810069624815
$ws.Range("B3").Value = 1730.2687525610861
$ws.Range("C3").Value = 1581.0506561184002
$ws.Range("A4").Value = 2494.283679208394
$ws.Range("B4").Value = 1874.7448644803062
$ws.Range("C4").Value = 1887.29228583497
$ws.Range("A5").Value = 2689.7593451447096
$ws.Range("B5").Value = 1755.0975951550595
$ws.Range("C5").Value = 1693.0873374066441
$ws.Range("A6").Value = 2478.7607831112336
$ws.Range("B6").Value = 1810.1830437242672
$ws.Range("C6").Value = 1877.1317007992568
$ws.Range("A7").Value = 2397.2088824637735
$ws.Range("B7").Value = 1847.2459982300322
$ws.Range("C7").Value = 1661.2757459849156
$ws.Range("A8").Value = 2467.9969983094284
$ws.Range("B8").Value = 1930.0313340409593
$ws.Range("C8").Value = 1785.2178304670479
$ws.Range("A9").Value = 2643.8788032222651
$ws.Range("B9").Value = 1943.5009763779306
$ws.Range("C9").Value = 1635.5503062445323
$ws.Range("A10").Value = 2394.9729326343991
$ws.Range("B10").Value = 1505.7671213936892
$ws.Range("C10").Value = 1488.8336686548521
$ws.Range("A11").Value = 2165.3642525000978
$ws.Range("B11").Value = 1568.9743546101593
$ws.Range("C11").Value = 1419.0680867570668
$ws.Range("A12").Value = 2864.4907222280667
$ws.Range("B12").Value = 2225.4880847403456
$ws.Range("C12").Value = 1936.9306010450307
$ws.Range("A13").Value = 2528.1335510759864
$ws.Range("B13").Value = 1931.1250729409865
$ws.Range("C13").Value = 1747.3103767153727
$ws.Range("A14").Value = 2609.7077692878297
$ws.Range("B14").Value = 2013.3383535139442
$ws.Range("C14").Value = 1762.6716978165794
$ws.Range("A15").Value = 2531.3178653380537
$ws.Range("B15").Value = 1976.917005535139
$ws.Range("C15").Value = 1827.4929104745929
$ws.Range("A16").Value = 2592.3756417833501
$ws.Range("B16").Value = 1796.2391444922016
$ws.Range("C16").Value = 1555.1927415373646
